# "Added Notification when plugged in and on alert"
# Phase2 gets a new leading outline-numbering column (1, 1.1, 1.2, 2, 2.1, ...)
# and a new item 4 / 4.1 describing the "Memory" / SharedPreferences requirement.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase2")

# Shift the whole sheet one column to the right to make room for the new
# outline-number column in column A.
$ws.Columns("A").Insert()

# Outline numbers for the existing requirement blocks.
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 1.1
$ws.Range("A3").Value = 1.2
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 2.1
$ws.Range("A10").Value = 3
$ws.Range("A11").Value = 3.1

# New requirement block 4: remember the entered value.
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "Memory"
$ws.Range("C15").Value = "Store entered value in SharedPreferences"
$ws.Range("A16").Value = 4.1

# Switch the page to portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection where Excel would land after typing the last entry.
[void]$ws.Range("A17").Select()
